$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 04:51:08"

# --- Sheet "zh-cn" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-20 04:50:55"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a27fd57a7621bb30a568f57da9f6708170871c2c/e2e/c91ca27e-9e4a-4fdf-9c2d-68d51cb06af5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b72c52e8cda2facffd69644016a2abc185497662/e2e/c91ca27e-9e4a-4fdf-9c2d-68d51cb06af5.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- Sheet "de-de" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-20 04:51:08"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a27fd57a7621bb30a568f57da9f6708170871c2c/e2e/c91ca27e-9e4a-4fdf-9c2d-68d51cb06af5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b72c52e8cda2facffd69644016a2abc185497662/e2e/c91ca27e-9e4a-4fdf-9c2d-68d51cb06af5.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
